$wb = $excel.ActiveWorkbook

# The file "c4c46589-1870-45fe-aa05-f073a8d8b456.md" has been handed off again.
# Update its status from "Handed back: in sync with en-US" to "Ready for handoff"
# on the Overview sheet, and update the per-language sheets with the new
# "Ready for handoff" status plus a refreshed "Latest Handoff Datetime".

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-02-18 03:40:44"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-02-18 03:40:57"
